$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the AutoFilter on the "DB_CourseNames__2" table: filter column "Kurs"
# down to a single course ("MAN 655  Corporate Strategy: Managing Business Groups")
# instead of the two previous "OPM 504 Transportation ..." entries.
# This also updates which rows are hidden (262 becomes visible, 321/322 become hidden).
$lo = $ws.ListObjects.Item(1)
[void]$lo.Range.AutoFilter(1, @("MAN 655  Corporate Strategy: Managing Business Groups"), 7)

# --- Re-apply the left/top-center alignment on the "Column1" notes cells that were
# carrying the redundant duplicate style (so they collapse back onto the single
# shared style instead of the stray duplicate with applyFill).
$noteCells = @("E8","E19","E85","E89","E115","E213","E214","E262","E320","E321","E328")
foreach ($addr in $noteCells) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4131
    $c.VerticalAlignment = -4108
}

# --- Move the active selection to E262 (the row that is now shown again).
[void]$ws.Range("E262").Select()
